$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.795.62"
$ws.Range("E2").Value = "  +4.76%  "

$ws.Range("D3").Value = "2.251.28"
$ws.Range("E3").Value = "  +3.97%  "

$ws.Range("E4").Value = "  +0.11%  "

$ws.Range("D5").Value = "'249.18"
$ws.Range("E5").Value = "  +1.01%  "

$ws.Range("D6").Value = "'0.633"
$ws.Range("E6").Value = "  +2.82%  "

$ws.Range("D7").Value = "'70.38"
$ws.Range("E7").Value = "  +6.16%  "

$ws.Range("E8").Value = "  -0.07%  "

$ws.Range("D9").Value = "'0.665"
$ws.Range("E9").Value = "  +17.29%  "

$ws.Range("D10").Value = "'39.06"
$ws.Range("E10").Value = "  +9.77%  "

$ws.Range("D11").Value = "'59.35"
$ws.Range("E11").Value = "  +2.00%  "

$ws.Range("D12").Value = "'0.0964"
$ws.Range("E12").Value = "  +4.34%  "

$ws.Range("E13").Value = "  +8.57%  "

$ws.Range("D14").Value = "'0.104"
$ws.Range("E14").Value = "  +1.34%  "

$ws.Range("D15").Value = "2.578.84"
$ws.Range("E15").Value = "  +3.76%  "

$ws.Range("D16").Value = "'14.79"
$ws.Range("E16").Value = "  +3.90%  "

$ws.Range("D17").Value = "'0.880"
$ws.Range("E17").Value = "  +2.38%  "

$ws.Range("D18").Value = "2.248.64"
$ws.Range("E18").Value = "  +4.83%  "

$ws.Range("D19").Value = "42.707.74"
$ws.Range("E19").Value = "  +4.68%  "

$ws.Range("D20").Value = "0.0₃0988"
$ws.Range("E20").Value = "  +5.45%  "

$ws.Range("E21").Value = "  +3.38%  "

$ws.Range("D22").Value = "'72.91"
$ws.Range("E22").Value = "  +2.38%  "

$ws.Range("D23").Value = "'235.07"
$ws.Range("E23").Value = "  +2.94%  "

$ws.Range("E24").Value = "  -3.51%  "

$ws.Range("E25").Value = "  +6.28%  "

$ws.Range("B26").Value = "Dai"
$ws.Range("C26").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  +0.06%  "

$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").Value = "'11.44"
$ws.Range("E27").Value = "  -2.26%  "

$ws.Range("D28").Value = "'2.41"
$ws.Range("E28").Value = "  -0.25%  "

$ws.Range("D29").Value = "'3.65"
$ws.Range("E29").Value = "  -1.57%  "

$ws.Range("D30").Value = "'2.11"
$ws.Range("E30").Value = "  +5.76%  "

$ws.Range("D31").Value = "'167.54"
$ws.Range("E31").Value = "  -0.62%  "

$ws.Range("E32").Value = "  +3.52%  "

$ws.Range("E33").Value = "  +14.55%  "

$ws.Range("D34").Value = "'0.125"
$ws.Range("E34").Value = "  +4.77%  "

$ws.Range("D35").Value = "'0.0800"
$ws.Range("E35").Value = "  +8.33%  "

$ws.Range("D36").Value = "'31.41"
$ws.Range("E36").Value = "  +27.25%  "

$ws.Range("D38").Value = "'4.45"
$ws.Range("E38").Value = "  +12.80%  "

$ws.Range("E39").Value = "  +3.57%  "

$ws.Range("D40").Value = "'0.0322"
$ws.Range("E40").Value = "  +8.39%  "

$ws.Range("D41").Value = "'2.30"
$ws.Range("E41").Value = "  +5.75%  "

$ws.Range("E42").Value = "  +7.64%  "

$ws.Range("D43").Value = "'5.78"
$ws.Range("E43").Value = "  +6.37%  "

$ws.Range("D44").Value = "'62.24"
$ws.Range("E44").Value = "  +3.11%  "

$ws.Range("E45").Value = "  +5.46%  "

$ws.Range("D46").Value = "'8.98"
$ws.Range("E46").Value = "  +6.07%  "

$ws.Range("E47").Value = "  +1.21%  "

$ws.Range("E48").Value = "  +2.96%  "

$ws.Range("E49").Value = "  -0.17%  "

$ws.Range("E50").Value = "  +2.19%  "

$ws.Range("E51").Value = "  +3.76%  "
